# Femacal de La Calera - Achicoria: weekly data update.
# A new observation (fecha serial 44460) is inserted as a new data row
# right before the existing row for fecha 44165 (current row 79),
# pushing that row and every row below it down by one. The last
# existing row simply moves down to become the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 79; Excel shifts rows 79..123 down to 80..124.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly observation.
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44460
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112010
$ws.Range("G79").Value = "Achicoria"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 60
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 6000
$ws.Range("M79").Value = 6000
$ws.Range("N79").Value = "$/caja 16 unidades"
$ws.Range("O79").Value = "Provincia de Quillota"
$ws.Range("P79").Value = 375
$ws.Range("Q79").Value = 16
$ws.Range("R79").Value = "Hortaliza"
